# Rename quiz sheets to numeric codes (e.g., 1.1)
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("1.1_Intro_to_Vectors").Name = "1.1"
$wb.Worksheets.Item("1.2_Vector_Addition").Name = "1.2"
$wb.Worksheets.Item("1.3_Dot_Product").Name = "1.3"
$wb.Worksheets.Item("1.4_Vectors_in_AI").Name = "1.4"

# Update the active selection on the "1.4" sheet (was the active tab) to C19.
$ws4 = $wb.Worksheets.Item("1.4")
$ws4.Activate()
$ws4.Range("C19").Select()
